# Insert a new data row before the current row 16 (shifting old rows 16-24
# down to 17-25), then populate the new row 16 with the latest weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16:24 down by one to make room for the new weekly record.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16.
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44455
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 100112026
$ws.Cells.Item(16, 7).Value = "Haba"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 10
$ws.Cells.Item(16, 11).Value = 13000
$ws.Cells.Item(16, 12).Value = 13000
$ws.Cells.Item(16, 13).Value = 13000
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 16).Value = 520
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"
